$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, shifting existing rows 9-22 down to 10-23.
$ws.Rows(9).Insert()

# Populate the newly inserted row 9 with the new weekly entry.
$ws.Cells.Item(9, 1).Value = 5
$ws.Cells.Item(9, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(9, 3).Value = "Maule"
$ws.Cells.Item(9, 4).Value = 44413
$ws.Cells.Item(9, 5).Value = 7
$ws.Cells.Item(9, 6).Value = 100112013
$ws.Cells.Item(9, 7).Value = "Alcachofa"
$ws.Cells.Item(9, 8).Value = "Española"
$ws.Cells.Item(9, 9).Value = "Segunda"
$ws.Cells.Item(9, 10).Value = 300
$ws.Cells.Item(9, 11).Value = 14000
$ws.Cells.Item(9, 12).Value = 14000
$ws.Cells.Item(9, 13).Value = 14000
$ws.Cells.Item(9, 14).Value = '$/caja 40 unidades'
$ws.Cells.Item(9, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(9, 16).Value = 350
$ws.Cells.Item(9, 17).Value = 40
$ws.Cells.Item(9, 18).Value = "Hortaliza"
